# TMTT0022150_VerificationOfCoverageCompanyDashboardNewUIAndFunctionality - 1 May 2024
#
# The "Users" sheet stores a single test-user display name in A2; the
# commit swaps the placeholder name for the real author, and leaves the
# sheet's selection sitting on F4 (as it was when the workbook was last
# saved by Excel).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

# A2: "James Craven" -> "Sahil Mittal"
$ws.Range("A2").Value = "Sahil Mittal"

# Leave the active selection on F4 (was D7).
$ws.Activate()
$ws.Range("F4").Select()
